$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell T1 ("5-jul"), same format as existing date headers (text, numFmtId 49)
$ws.Cells.Item(1, 20).NumberFormat = "@"
$ws.Cells.Item(1, 20).Value = "5-jul"

# New numeric values for column T, rows 2-11, matching the format of column S
# (integer number format, centered horizontal alignment)
$tValues = @(13, 23, 11, 12, 13, 17, 11, 21, 23, 6)

for ($i = 0; $i -lt $tValues.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 20)
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
    $cell.Value = $tValues[$i]
}

# Update the selected cell to match the new saved state
$ws.Range("F16").Select()
